# RP3_APT_TxOut_2024_Jan_Dec.xlsx - March release update + 2024 post ops
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TxOut_APT")

# Updated "last refreshed" date in the header block (B2)
$ws.Range("B2").Value = 45758.0

# Block 1 (rows 6-17): column E used to hold the shared formula F/D.
# It is now a plain value rounded to 2 decimals; column F is rounded to
# the nearest whole number (some rows' F values are unchanged and are
# skipped below).
$ws.Range("E6").Value  = 2.08
$ws.Range("F6").Value  = 198153.0
$ws.Range("E7").Value  = 1.74
$ws.Range("F7").Value  = 155226.0
$ws.Range("E8").Value  = 3.24
$ws.Range("F8").Value  = 691174.0
$ws.Range("E9").Value  = 1.35
$ws.Range("F9").Value  = 81236.0
$ws.Range("E10").Value = 1.19
$ws.Range("F10").Value = 66560.0
$ws.Range("E11").Value = 1.89
$ws.Range("F11").Value = 144234.0
$ws.Range("E12").Value = 2.9
$ws.Range("E13").Value = 1.57
$ws.Range("F13").Value = 69997.0
$ws.Range("E14").Value = 2.9
$ws.Range("F14").Value = 219342.0
$ws.Range("E15").Value = 3.56
$ws.Range("F15").Value = 851211.0
$ws.Range("E16").Value = 3.65
$ws.Range("F16").Value = 436518.0
$ws.Range("E17").Value = 2.53
$ws.Range("F17").Value = 288353.0

# Block 2 (rows 19-29)
$ws.Range("E19").Value = 3.44
$ws.Range("E20").Value = 2.78
$ws.Range("F20").Value = 251627.0
$ws.Range("E21").Value = 1.81
$ws.Range("F21").Value = 166204.0
$ws.Range("E22").Value = 2.32
$ws.Range("F22").Value = 156090.0
$ws.Range("E23").Value = 1.67
$ws.Range("F23").Value = 95785.0
$ws.Range("E24").Value = 3.18
$ws.Range("F24").Value = 544220.0
$ws.Range("E25").Value = 3.3
$ws.Range("F25").Value = 653284.0
$ws.Range("E26").Value = 2.26
$ws.Range("F26").Value = 187035.0
$ws.Range("E27").Value = 2.48
$ws.Range("F27").Value = 296398.0
$ws.Range("E28").Value = 0.84
$ws.Range("F28").Value = 27790.0
$ws.Range("E29").Value = 0.96

# Block 3 (rows 31-47)
$ws.Range("E31").Value = 2.08
$ws.Range("F31").Value = 154559.0
$ws.Range("E32").Value = 3.75
$ws.Range("F32").Value = 810385.0
$ws.Range("E33").Value = 1.91
$ws.Range("E34").Value = 2.95
$ws.Range("F34").Value = 384413.0
$ws.Range("E35").Value = 1.31
$ws.Range("F35").Value = 80958.0
$ws.Range("E36").Value = 3.54
$ws.Range("F36").Value = 361252.0
$ws.Range("E37").Value = 2.09
$ws.Range("E38").Value = 3.71
$ws.Range("F38").Value = 214464.0
$ws.Range("E39").Value = 1.67
$ws.Range("E40").Value = 7.28
$ws.Range("E41").Value = 2.17
$ws.Range("E42").Value = 2.7
$ws.Range("F42").Value = 329122.0
$ws.Range("E43").Value = 1.86
$ws.Range("F43").Value = 99462.0
$ws.Range("E44").Value = 4.32
$ws.Range("F44").Value = 470453.0
$ws.Range("E45").Value = 2.24
$ws.Range("F45").Value = 121772.0
$ws.Range("E46").Value = 2.68
$ws.Range("F46").Value = 235090.0
$ws.Range("E47").Value = 3.04
$ws.Range("F47").Value = 387227.0
